# Update test fixture for RT (related table) importer: add a new "Group"
# column (F) to the wide related-table sheet. Column F mirrors the
# geom_code value in column A for rows whose code is "A" or "B", and is
# left blank for rows whose code is "C".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the gridlines display the same as before touching the sheet (the
# runtime otherwise flips the (default-valued) attribute when it rewrites
# the view on save).
$excel.ActiveWindow.DisplayGridlines = $true

# New header for column F.
$ws.Range("F1").Value = "Group"

# Widen column E (DateFormat) - the new data in column F uses the
# sheet's default column width.
$ws.Columns.Item(5).ColumnWidth = 19.4

# Mirror column A (geom_code) into the new column F, row by row, except
# rows where geom_code is "C" which stay blank (but still present as an
# empty, formatted cell).
for ($r = 2; $r -le 10; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    $cell = $ws.Cells.Item($r, 6)
    if ($code -eq "C") {
        $cell.NumberFormat = "General"
    } else {
        $cell.Value = $code
    }
}

# Best-effort page setup touch-up (matches the fitToWidth/fitToHeight
# that stay set in the target file).
$ps = $ws.PageSetup
$ps.FitToPagesWide = 1
$ps.FitToPagesTall = 1

# Leave the cursor on the last edited cell, matching the saved selection.
$ws.Range("F10").Select()
